$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-05-28T14:08:24"
$ws.Range("S4").Value = 59.41
$ws.Range("T4").Value = 32.27
$ws.Range("U4").Value = 48.85
$ws.Range("W4").Value = 31.33
$ws.Range("X4").Value = 27.71
$ws.Range("Y4").Value = 14.78
$ws.Range("R6").Value = -0.31
$ws.Range("S6").Value = -0.71
$ws.Range("T6").Value = -0.39
$ws.Range("U6").Value = -0.64
$ws.Range("W6").Value = -0.6
$ws.Range("X6").Value = -0.19
$ws.Range("Y6").Value = 0
$ws.Range("R9").Value = 33.13
$ws.Range("S9").Value = 56.72
$ws.Range("T9").Value = 31.28
$ws.Range("U9").Value = 47.68
$ws.Range("V9").Value = 32.88
$ws.Range("W9").Value = 30.37
$ws.Range("X9").Value = 27.57
$ws.Range("Y9").Value = 15.1
$ws.Range("Z9").Value = 15.02
$ws.Range("R11").Value = -2.19
$ws.Range("S11").Value = -3.4
$ws.Range("T11").Value = -1.38
$ws.Range("U11").Value = -1.81
$ws.Range("V11").Value = -1.28
$ws.Range("W11").Value = -1.55
$ws.Range("X11").Value = -0.33
$ws.Range("Y11").Value = 0.32
$ws.Range("Z11").Value = 0.36
$ws.Range("R14").Value = 14.35
$ws.Range("S14").Value = 57.49
$ws.Range("T14").Value = 31.25
$ws.Range("U14").Value = 47.68
$ws.Range("V14").Value = 32.85
$ws.Range("W14").Value = 30.37
$ws.Range("X14").Value = 27.54
$ws.Range("Z14").Value = 15.02
$ws.Range("R15").Value = -18.74
$ws.Range("S15").Value = 0.77
$ws.Range("T15").Value = 0
$ws.Range("S16").Value = -3.4
$ws.Range("T16").Value = -1.41
$ws.Range("U16").Value = -1.81
$ws.Range("V16").Value = -1.31
$ws.Range("W16").Value = -1.55
$ws.Range("X16").Value = -0.36
$ws.Range("Z16").Value = 0.36
$ws.Range("S19").Value = 59.47
$ws.Range("T19").Value = 32.34
$ws.Range("U19").Value = 49
$ws.Range("W19").Value = 31.39
$ws.Range("X19").Value = 27.93
$ws.Range("Z19").Value = 14.89
$ws.Range("R21").Value = -0.28
$ws.Range("S21").Value = -0.65
$ws.Range("U21").Value = -0.49
$ws.Range("W21").Value = -0.53
$ws.Range("Z21").Value = 0.22
$ws.Range("S24").Value = 59.47
$ws.Range("T24").Value = 32.34
$ws.Range("U24").Value = 49
$ws.Range("W24").Value = 31.39
$ws.Range("X24").Value = 27.93
$ws.Range("Z24").Value = 14.89
$ws.Range("R26").Value = -0.28
$ws.Range("S26").Value = -0.65
$ws.Range("U26").Value = -0.49
$ws.Range("W26").Value = -0.53
$ws.Range("Z26").Value = 0.22
$ws.Range("S29").Value = 59.76
$ws.Range("T29").Value = 32.56
$ws.Range("U29").Value = 49.34
$ws.Range("V29").Value = 33.93
$ws.Range("W29").Value = 31.61
$ws.Range("X29").Value = 28.24
$ws.Range("Z29").Value = 15.1
$ws.Range("R31").Value = -0.11
$ws.Range("S31").Value = -0.36
$ws.Range("U31").Value = -0.15
$ws.Range("V31").Value = -0.24
$ws.Range("W31").Value = -0.32
$ws.Range("X31").Value = 0.34
$ws.Range("Z31").Value = 0.44
$ws.Range("R34").Value = 13.02
$ws.Range("S34").Value = 55.42
$ws.Range("T34").Value = 30.35
$ws.Range("U34").Value = 46.64
$ws.Range("V34").Value = 32.11
$ws.Range("W34").Value = 29.42
$ws.Range("X34").Value = 27.27
$ws.Range("R35").Value = -18.74
$ws.Range("S35").Value = 0.77
$ws.Range("T35").Value = 0
$ws.Range("R36").Value = -3.56
$ws.Range("S36").Value = -5.47
$ws.Range("T36").Value = -2.31
$ws.Range("U36").Value = -2.85
$ws.Range("V36").Value = -2.06
$ws.Range("W36").Value = -2.5
$ws.Range("X36").Value = -0.63
$ws.Range("S39").Value = 59.41
$ws.Range("T39").Value = 32.27
$ws.Range("U39").Value = 48.85
$ws.Range("W39").Value = 31.33
$ws.Range("X39").Value = 27.71
$ws.Range("Y39").Value = 14.78
$ws.Range("R41").Value = -0.31
$ws.Range("S41").Value = -0.71
$ws.Range("T41").Value = -0.39
$ws.Range("U41").Value = -0.64
$ws.Range("W41").Value = -0.6
$ws.Range("X41").Value = -0.19
$ws.Range("Y41").Value = 0
$ws.Range("S44").Value = 60.79
$ws.Range("T44").Value = 33.06
$ws.Range("U44").Value = 50.09
$ws.Range("W44").Value = 32.18
$ws.Range("X44").Value = 28.01
$ws.Range("R46").Value = 0.32
$ws.Range("S46").Value = 0.67
$ws.Range("U46").Value = 0.6
$ws.Range("W46").Value = 0.26
$ws.Range("X46").Value = 0.11
$ws.Range("R49").Value = 31.64
$ws.Range("S49").Value = 54.75
$ws.Range("T49").Value = 28.2
$ws.Range("U49").Value = 42.44
$ws.Range("V49").Value = 29.4
$ws.Range("W49").Value = 27.4
$ws.Range("X49").Value = 23.65
$ws.Range("R51").Value = -3.67
$ws.Range("S51").Value = -5.37
$ws.Range("U51").Value = -7.05
$ws.Range("V51").Value = -4.76
$ws.Range("W51").Value = -4.52
$ws.Range("X51").Value = -4.26
$ws.Range("R54").Value = 29.53
$ws.Range("S54").Value = 50.56
$ws.Range("T54").Value = 27.13
$ws.Range("U54").Value = 41.27
$ws.Range("V54").Value = 29.28
$ws.Range("W54").Value = 26.62
$ws.Range("X54").Value = 23.99
$ws.Range("Z54").Value = 12.56
$ws.Range("R56").Value = -5.79
$ws.Range("S56").Value = -9.56
$ws.Range("T56").Value = -5.53
$ws.Range("U56").Value = -8.210000000000001
$ws.Range("V56").Value = -4.89
$ws.Range("W56").Value = -5.3
$ws.Range("X56").Value = -3.91
$ws.Range("Z56").Value = -2.1
$ws.Range("R59").Value = 36.67
$ws.Range("S59").Value = 62.5
$ws.Range("T59").Value = 33.95
$ws.Range("U59").Value = 51.44
$ws.Range("W59").Value = 33.08
$ws.Range("X59").Value = 28.74
$ws.Range("S61").Value = 2.37
$ws.Range("U61").Value = 1.95
$ws.Range("W61").Value = 1.16
$ws.Range("X61").Value = 0.83
$ws.Range("R64").Value = 37.25
$ws.Range("S64").Value = 63.49
$ws.Range("T64").Value = 34.49
$ws.Range("U64").Value = 52.26
$ws.Range("W64").Value = 33.67
$ws.Range("X64").Value = 29.16
$ws.Range("R66").Value = 1.94
$ws.Range("S66").Value = 3.36
$ws.Range("U66").Value = 2.77
$ws.Range("W66").Value = 1.75
$ws.Range("X66").Value = 1.25
$ws.Range("R69").Value = 37.69
$ws.Range("S69").Value = 64.16
$ws.Range("T69").Value = 34.86
$ws.Range("U69").Value = 52.82
$ws.Range("W69").Value = 34
$ws.Range("X69").Value = 29.43
$ws.Range("R71").Value = 2.37
$ws.Range("S71").Value = 4.04
$ws.Range("U71").Value = 3.33
$ws.Range("W71").Value = 2.07
$ws.Range("X71").Value = 1.53
$ws.Range("R74").Value = 35.96
$ws.Range("S74").Value = 61.35
$ws.Range("T74").Value = 33.26
$ws.Range("U74").Value = 50.4
$ws.Range("W74").Value = 32.44
$ws.Range("X74").Value = 28.13
$ws.Range("S76").Value = 1.23
$ws.Range("U76").Value = 0.91
$ws.Range("W76").Value = 0.52
$ws.Range("X76").Value = 0.23
$ws.Range("R79").Value = 35.31
$ws.Range("S79").Value = 60.12
$ws.Range("T79").Value = 32.66
$ws.Range("U79").Value = 49.49
$ws.Range("W79").Value = 31.92
$ws.Range("X79").Value = 27.9
$ws.Range("R84").Value = 29.53
$ws.Range("S84").Value = 50.61
$ws.Range("T84").Value = 27.13
$ws.Range("U84").Value = 41.31
$ws.Range("V84").Value = 29.28
$ws.Range("W84").Value = 26.65
$ws.Range("X84").Value = 24.01
$ws.Range("Z84").Value = 12.56
$ws.Range("R86").Value = -5.79
$ws.Range("S86").Value = -9.51
$ws.Range("T86").Value = -5.53
$ws.Range("U86").Value = -8.18
$ws.Range("V86").Value = -4.89
$ws.Range("W86").Value = -5.28
$ws.Range("X86").Value = -3.89
$ws.Range("Z86").Value = -2.1
$ws.Range("S89").Value = 59.76
$ws.Range("T89").Value = 32.56
$ws.Range("U89").Value = 49.34
$ws.Range("V89").Value = 33.93
$ws.Range("W89").Value = 31.61
$ws.Range("X89").Value = 28.27
$ws.Range("R91").Value = -0.11
$ws.Range("S91").Value = -0.36
$ws.Range("U91").Value = -0.15
$ws.Range("V91").Value = -0.24
$ws.Range("W91").Value = -0.32
$ws.Range("X91").Value = 0.37
